$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting existing data (D:K) to (E:L)
$ws.Columns("D").Insert()

# Copy number formats/styles from column E (the old column D, now shifted) into new column D
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)

# Populate the new column D with the latest period values
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 2448400
$ws.Range("D9").Value = 1063900
$ws.Range("D10").Value = 1384400
$ws.Range("D12").Value = 106200
$ws.Range("D14").Value = 79200
$ws.Range("D17").Value = 2126700
$ws.Range("D18").Value = 321700
$ws.Range("D20").Value = 900
$ws.Range("D21").Value = 532600
$ws.Range("D22").Value = 103000
$ws.Range("D23").Value = 219600
$ws.Range("D24").Value = 23200
$ws.Range("D26").Value = 196400
$ws.Range("D27").Value = 196400
$ws.Range("D29").Value = 4400
$ws.Range("D32").Value = -900
$ws.Range("D33").Value = 200800
$ws.Range("D35").Value = 200800
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 357200
$ws.Range("D43").Value = 366300
$ws.Range("D44").Value = 427800
$ws.Range("D45").Value = 84900
$ws.Range("D46").Value = 1236200
$ws.Range("D47").Value = "NA"
$ws.Range("D48").Value = 432800
$ws.Range("D49").Value = 4571600
$ws.Range("D52").Value = 37400
$ws.Range("D54").Value = 6278000
$ws.Range("D57").Value = 106700
$ws.Range("D58").Value = 86600
$ws.Range("D59").Value = 389100
$ws.Range("D60").Value = 582500
$ws.Range("D61").Value = 2072200
$ws.Range("D62").Value = 1083400
$ws.Range("D66").Value = 3738000
$ws.Range("D72").Value = 2427600
$ws.Range("D76").Value = 2540000
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 200800
$ws.Range("D83").Value = 210000
$ws.Range("D89").Value = 437400
$ws.Range("D91").Value = -80800
$ws.Range("D94").Value = -196400
$ws.Range("D96").Value = -62200
$ws.Range("D100").Value = -206400
$ws.Range("D101").Value = -10900
$ws.Range("D102").Value = 23600
$ws.Range("D13").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("D28").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D34").Value = 0
$ws.Range("D42").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D53").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D77").Value = 0
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
